# Auto-generated script applying the Omega_Profits cell-value updates
# described by the commit diff. Each sheet is addressed by its tab name;
# cells are updated via Range(...).Value assignments, matching the exact
# before -> after numeric deltas from the diff. Two cells are net-new
# (N134 on LTW, N126 on WVR) and three are net-removed (M63/M66 on CUL,
# M117 on CUL) -- those use ClearContents() to delete the cell entirely.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1200.1818
$ws.Range("I28").Value = 680.5
$ws.Range("J28").Value = 2586
$ws.Range("K28").Value = 680.5
$ws.Range("L28").Value = 2586
$ws.Range("M28").Value = -195.5
$ws.Range("N28").Value = -3556
$ws.Range("H42").Value = 2223.1428
$ws.Range("I42").Value = 638.75
$ws.Range("J42").Value = 4335.6665
$ws.Range("K42").Value = 1916.25
$ws.Range("L42").Value = 13006.9995
$ws.Range("M42").Value = -1686.25
$ws.Range("N42").Value = -13466.9995
$ws.Range("H43").Value = 6087.6665
$ws.Range("I43").Value = 11959.2
$ws.Range("J43").Value = 3151.9
$ws.Range("K43").Value = 11959.2
$ws.Range("L43").Value = 3151.9
$ws.Range("M43").Value = -11890.2
$ws.Range("N43").Value = -3289.9
$ws.Range("H76").Value = 8115.44
$ws.Range("I76").Value = 8718.846
$ws.Range("J76").Value = 7461.75
$ws.Range("K76").Value = 8718.846
$ws.Range("L76").Value = 7461.75
$ws.Range("M76").Value = -8403.846
$ws.Range("N76").Value = -8091.75
$ws.Range("H79").Value = 8115.44
$ws.Range("I79").Value = 8718.846
$ws.Range("J79").Value = 7461.75
$ws.Range("K79").Value = 8718.846
$ws.Range("L79").Value = 7461.75
$ws.Range("M79").Value = -7626.846
$ws.Range("N79").Value = -9645.75
$ws.Range("H86").Value = 3522.7144
$ws.Range("I86").Value = 3931.6
$ws.Range("K86").Value = 3931.6
$ws.Range("M86").Value = -2808.6
$ws.Range("H89").Value = 3522.7144
$ws.Range("I89").Value = 3931.6
$ws.Range("K89").Value = 19658
$ws.Range("M89").Value = -14042
$ws.Range("H138").Value = 3658.027
$ws.Range("J138").Value = 4161.6
$ws.Range("L138").Value = 12484.8
$ws.Range("N138").Value = -22764.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9981.575000000001
$ws.Range("I32").Value = 940.2759
$ws.Range("K32").Value = 940.2759
$ws.Range("M32").Value = -653.2759
$ws.Range("H74").Value = 2609.4375
$ws.Range("I74").Value = 1800.5454
$ws.Range("K74").Value = 1800.5454
$ws.Range("M74").Value = -926.5454
$ws.Range("H77").Value = 2609.4375
$ws.Range("I77").Value = 1800.5454
$ws.Range("K77").Value = 9002.726999999999
$ws.Range("M77").Value = -4634.726999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 1163.5454
$ws.Range("I19").Value = 1382.1111
$ws.Range("K19").Value = 1382.1111
$ws.Range("M19").Value = -1209.1111
$ws.Range("H80").Value = 1042.2174
$ws.Range("I80").Value = 1325.5385
$ws.Range("J80").Value = 673.9
$ws.Range("K80").Value = 1325.5385
$ws.Range("L80").Value = 673.9
$ws.Range("M80").Value = -327.5385000000001
$ws.Range("N80").Value = -2669.9
$ws.Range("H83").Value = 1042.2174
$ws.Range("I83").Value = 1325.5385
$ws.Range("J83").Value = 673.9
$ws.Range("K83").Value = 6627.692500000001
$ws.Range("L83").Value = 3369.5
$ws.Range("M83").Value = -1635.692500000001
$ws.Range("N83").Value = -13353.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1919.2
$ws.Range("I22").Value = 1399.25
$ws.Range("K22").Value = 1399.25
$ws.Range("M22").Value = -1049.25
$ws.Range("H99").Value = 3090778
$ws.Range("I99").Value = 3972740.2
$ws.Range("J99").Value = 3909.5
$ws.Range("K99").Value = 3972740.2
$ws.Range("L99").Value = 3909.5
$ws.Range("M99").Value = -3971242.2
$ws.Range("N99").Value = -6905.5
$ws.Range("H126").Value = 3090778
$ws.Range("I126").Value = 3972740.2
$ws.Range("J126").Value = 3909.5
$ws.Range("K126").Value = 11918220.6
$ws.Range("L126").Value = 11728.5
$ws.Range("M126").Value = -11915750.6
$ws.Range("N126").Value = -16668.5
$ws.Range("H132").Value = 5558.857
$ws.Range("I132").Value = 4382.9
$ws.Range("K132").Value = 13148.7
$ws.Range("M132").Value = -10618.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 125007864
$ws.Range("I8").Value = 125007864
$ws.Range("K8").Value = 375023592
$ws.Range("M8").Value = -375023453
$ws.Range("H57").Value = 15667.333
$ws.Range("J57").Value = 18199.8
$ws.Range("L57").Value = 54599.39999999999
$ws.Range("N57").Value = -55717.39999999999
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()  # was -9751
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()  # was -27756
$ws.Range("H75").Value = 2841.1667
$ws.Range("I75").Value = 1789
$ws.Range("K75").Value = 5367
$ws.Range("M75").Value = -4369
$ws.Range("H78").Value = 2841.1667
$ws.Range("I78").Value = 1789
$ws.Range("K78").Value = 16101
$ws.Range("M78").Value = -11109
$ws.Range("H97").Value = 456.1111
$ws.Range("I97").Value = 388.25
$ws.Range("K97").Value = 1164.75
$ws.Range("M97").Value = -668.75
$ws.Range("H98").Value = 1411.25
$ws.Range("I98").Value = 1111.625
$ws.Range("J98").Value = 1710.875
$ws.Range("K98").Value = 3334.875
$ws.Range("L98").Value = 5132.625
$ws.Range("M98").Value = -1836.875
$ws.Range("N98").Value = -8128.625
$ws.Range("H117").Value = 5255
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 5255
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 15765
$ws.Range("M117").ClearContents()  # was 2575
$ws.Range("N117").Value = -22649
$ws.Range("H122").Value = 9087.25
$ws.Range("J122").Value = 5599.7144
$ws.Range("L122").Value = 50397.4296
$ws.Range("N122").Value = -55297.4296
$ws.Range("H136").Value = 11735.091
$ws.Range("I136").Value = 11723.5
$ws.Range("J136").Value = 11766
$ws.Range("K136").Value = 35170.5
$ws.Range("L136").Value = 35298
$ws.Range("M136").Value = -30070.5
$ws.Range("N136").Value = -45498
$ws.Range("H137").Value = 5088.846
$ws.Range("I137").Value = 4073.3333
$ws.Range("J137").Value = 7373.75
$ws.Range("K137").Value = 12219.9999
$ws.Range("L137").Value = 22121.25
$ws.Range("M137").Value = -7119.999899999999
$ws.Range("N137").Value = -32321.25
$ws.Range("H140").Value = 2569.2
$ws.Range("I140").Value = 2413.56
$ws.Range("K140").Value = 7240.68
$ws.Range("M140").Value = -2060.68

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 45924
$ws.Range("J106").Value = 45924
$ws.Range("L106").Value = 45924
$ws.Range("N106").Value = -48448

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7266.5
$ws.Range("I46").Value = 6708.125
$ws.Range("K46").Value = 6708.125
$ws.Range("M46").Value = -6520.125
$ws.Range("H93").Value = 779
$ws.Range("J93").Value = 197.5
$ws.Range("L93").Value = 197.5
$ws.Range("N93").Value = -2693.5
$ws.Range("H132").Value = 44987.65
$ws.Range("I132").Value = 53519.5
$ws.Range("J132").Value = 5172.3335
$ws.Range("K132").Value = 160558.5
$ws.Range("L132").Value = 15517.0005
$ws.Range("M132").Value = -158028.5
$ws.Range("N132").Value = -20577.0005
$ws.Range("H134").Value = 89999
$ws.Range("J134").Value = 89999
$ws.Range("L134").Value = 89999
$ws.Range("N134").Value = -100139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17103
$ws.Range("J45").Value = 17140.875
$ws.Range("L45").Value = 17140.875
$ws.Range("N45").Value = -18122.875
$ws.Range("H81").Value = 2998.8696
$ws.Range("J81").Value = 1713.2858
$ws.Range("L81").Value = 3426.5716
$ws.Range("N81").Value = -5548.5716
$ws.Range("H84").Value = 2998.8696
$ws.Range("J84").Value = 1713.2858
$ws.Range("L84").Value = 17132.858
$ws.Range("N84").Value = -27740.858
$ws.Range("H107").Value = 1674.875
$ws.Range("I107").Value = 1198.3125
$ws.Range("J107").Value = 2151.4375
$ws.Range("K107").Value = 3594.9375
$ws.Range("L107").Value = 6454.3125
$ws.Range("M107").Value = -1674.9375
$ws.Range("N107").Value = -10294.3125
$ws.Range("H126").Value = 2999.6667
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -19937
$ws.Range("H132").Value = 2985.7878
$ws.Range("I132").Value = 2987.7334
$ws.Range("K132").Value = 8963.200199999999
$ws.Range("M132").Value = -6433.200199999999
